$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.5898876404494382
$ws1.Range("C2").Value = 0.8478260869565217
$ws1.Range("D2").Value = 0.2191011235955056
$ws1.Range("E2").Value = 0.3482142857142857
$ws1.Range("F2").Value = 0.5898876404494382
$ws1.Range("G2").Value = 117
$ws1.Range("H2").Value = 21
$ws1.Range("I2").Value = 513
$ws1.Range("J2").Value = 417

# --- Sheet: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")

$ws2.Range("B2").Value = 0.5516129032258065
$ws2.Range("C2").Value = 0.9606741573033708
$ws2.Range("D2").Value = 0.7008196721311475

$ws2.Range("B3").Value = 0.8478260869565217
$ws2.Range("C3").Value = 0.2191011235955056
$ws2.Range("D3").Value = 0.3482142857142857
$ws2.Range("E3").Value = 534

$ws2.Range("B4").Value = 0.5898876404494382
$ws2.Range("C4").Value = 0.5898876404494382
$ws2.Range("D4").Value = 0.5898876404494382
$ws2.Range("E4").Value = 0.5898876404494382

$ws2.Range("B5").Value = 0.699719495091164
$ws2.Range("C5").Value = 0.5898876404494382
$ws2.Range("D5").Value = 0.5245169789227166
$ws2.Range("E5").Value = 1068

$ws2.Range("B6").Value = 0.6997194950911642
$ws2.Range("C6").Value = 0.5898876404494382
$ws2.Range("D6").Value = 0.5245169789227166
$ws2.Range("E6").Value = 1068

# --- Sheet: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

$ws3.Range("B2").Value = 513
$ws3.Range("C2").Value = 21

$ws3.Range("B3").Value = 417
$ws3.Range("C3").Value = 117
